# Fixing raw correlation table: blank out the stale/incorrect cells in
# the "sample-sample raw correlation" table on slide 2.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tbl = $shp.Table

# (row, col) pairs are 1-based and include the header row/column.
$targets = @(
    @(2, 4),  # mouse_fibrosis  x human_fibrosis  -> was 0.85
    @(4, 2),  # human_fibrosis  x mouse_fibrosis  -> was 0.87
    @(4, 3),  # human_fibrosis  x mouse_control   -> was 0.84
    @(5, 2),  # human_control   x mouse_fibrosis  -> was 0.86
    @(5, 3),  # human_control   x mouse_control   -> was 0.88
    @(5, 4)   # human_control   x human_fibrosis  -> was 0.87
)

foreach ($t in $targets) {
    $row = $t[0]
    $col = $t[1]
    $cell = $tbl.Cell($row, $col)
    $cell.Shape.TextFrame.TextRange.Text = ""
}
